$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($ws, $addr, $val) {
    $r = $ws.Range($addr)
    $r.NumberFormat = "@"
    $r.Value = $val
    $r.Style = "Normal"
}

Set-TextValue $ws "D2" "328.49"
Set-TextValue $ws "E2" "-0.76%"
Set-TextValue $ws "D3" "43.93"
Set-TextValue $ws "E3" "5.36%"
Set-TextValue $ws "E4" "-4.17%"
Set-TextValue $ws "D5" "0.08110"
Set-TextValue $ws "E5" "-2.89%"
Set-TextValue $ws "D6" "8.700"
Set-TextValue $ws "E6" "-1.08%"
Set-TextValue $ws "D7" "4.320"
Set-TextValue $ws "E7" "-3.20%"
Set-TextValue $ws "D8" "1.895"
Set-TextValue $ws "E8" "-4.73%"
Set-TextValue $ws "D9" "2.768"
Set-TextValue $ws "E9" "-4.87%"
Set-TextValue $ws "D10" "0.9453"
Set-TextValue $ws "E10" "2.18%"
Set-TextValue $ws "D11" "0.1183"
Set-TextValue $ws "E11" "-8.11%"
Set-TextValue $ws "D12" "0.1889"
Set-TextValue $ws "E12" "-4.44%"
Set-TextValue $ws "D13" "0.09617"
Set-TextValue $ws "E13" "1.96%"
Set-TextValue $ws "D14" "0.04221"
Set-TextValue $ws "E14" "9.21%"
Set-TextValue $ws "D15" "0.1071"
Set-TextValue $ws "E15" "0.98%"
Set-TextValue $ws "D16" "0.001285"
Set-TextValue $ws "E16" "-1.00%"
Set-TextValue $ws "D17" "0.005982"
Set-TextValue $ws "E17" "-1.98%"
Set-TextValue $ws "D18" "3.564"
Set-TextValue $ws "E18" "3.56%"
Set-TextValue $ws "E19" "-0.51%"
Set-TextValue $ws "D20" "8.700"
Set-TextValue $ws "E20" "2.48%"
Set-TextValue $ws "D21" "0.1361"
Set-TextValue $ws "E21" "-0.01%"
Set-TextValue $ws "D22" "0.2607"
Set-TextValue $ws "E22" "5.00%"
Set-TextValue $ws "D23" "0.04381"
Set-TextValue $ws "E23" "-0.71%"
Set-TextValue $ws "D24" "0.001243"
Set-TextValue $ws "E24" "-2.35%"
Set-TextValue $ws "D25" "0.004305"
Set-TextValue $ws "E25" "-1.70%"
Set-TextValue $ws "E26" "1.46%"
Set-TextValue $ws "D27" "0.0004017"
Set-TextValue $ws "E27" "31.92%"
Set-TextValue $ws "D39" "0.02718"
Set-TextValue $ws "E39" "-5.40%"
Set-TextValue $ws "D40" "0.05546"
Set-TextValue $ws "E40" "0.25%"
Set-TextValue $ws "D41" "0.007816"
Set-TextValue $ws "E41" "-1.48%"
Set-TextValue $ws "D42" "0.009757"
Set-TextValue $ws "E42" "5.03%"
Set-TextValue $ws "D43" "0.1397"
Set-TextValue $ws "E43" "-2.67%"
Set-TextValue $ws "E44" "-5.14%"
Set-TextValue $ws "D45" "0.01035"
Set-TextValue $ws "E45" "-6.84%"
Set-TextValue $ws "D46" "0.00007111"
Set-TextValue $ws "E46" "0.12%"
Set-TextValue $ws "D47" "0.00000000755"
Set-TextValue $ws "E47" "0.74%"
Set-TextValue $ws "D48" "0.003473"
Set-TextValue $ws "E48" "0.81%"
Set-TextValue $ws "D49" "0.002286"
Set-TextValue $ws "E49" "0.41%"
Set-TextValue $ws "E50" "0.74%"
Set-TextValue $ws "D51" "0.0002014"
Set-TextValue $ws "E51" "0.74%"
